$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "data-user-id"
$ws.Range("A2").Value = "TC-B06-001"
$ws.Range("B2").Value = 2
$ws.Range("A3").Value = "TC-B06-002,`nTC-B06-003,`nTC-B06-004"
$ws.Range("B3").Value = 6

# Font size changes on cells that already use Times New Roman (keeps family=1)
$ws.Range("A2").Font.Size = 11
$ws.Range("A3").Font.Size = 11

# Font name change on cells using the default Calibri font
$ws.Range("A1:B1").Font.Name = "Times New Roman"
$ws.Range("A1:B1").Font.Family = 1
$ws.Range("B2").Font.Name = "Times New Roman"
$ws.Range("B2").Font.Family = 1
$ws.Range("B3").Font.Name = "Times New Roman"
$ws.Range("B3").Font.Family = 1
